# Auto-generated edit script: updates crypto price/volume/hour data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'329.25"
$ws.Range("E2").Value = "'0.16%"
$ws.Range("G2").Value = "'18"
$ws.Range("D3").Value = "'44.48"
$ws.Range("E3").Value = "'0.30%"
$ws.Range("G3").Value = "'18"
$ws.Range("D4").Value = "'5.513"
$ws.Range("E4").Value = "'-1.23%"
$ws.Range("G4").Value = "'18"
$ws.Range("D5").Value = "'0.08121"
$ws.Range("E5").Value = "'0.36%"
$ws.Range("G5").Value = "'18"
$ws.Range("D6").Value = "'2.058"
$ws.Range("E6").Value = "'5.39%"
$ws.Range("G6").Value = "'18"
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").Value = "'4.434"
$ws.Range("E7").Value = "'2.82%"
$ws.Range("G7").Value = "'18"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.9677"
$ws.Range("E8").Value = "'1.57%"
$ws.Range("G8").Value = "'18"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "'0.1137"
$ws.Range("E9").Value = "'-3.28%"
$ws.Range("G9").Value = "'18"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1880"
$ws.Range("E10").Value = "'1.47%"
$ws.Range("G10").Value = "'18"
$ws.Range("B11").Value = "MCDex"
$ws.Range("C11").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D11").Value = "'10.14"
$ws.Range("E11").Value = "'-0.18%"
$ws.Range("G11").Value = "'18"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.09949"
$ws.Range("E12").Value = "'0.95%"
$ws.Range("G12").Value = "'18"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.04694"
$ws.Range("E13").Value = "'4.16%"
$ws.Range("G13").Value = "'18"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.1061"
$ws.Range("E14").Value = "'-0.70%"
$ws.Range("G14").Value = "'18"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001262"
$ws.Range("E15").Value = "'-1.51%"
$ws.Range("G15").Value = "'18"
$ws.Range("B16").Value = "CoinExToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D16").Value = "'0.04090"
$ws.Range("E16").Value = "'-2.71%"
$ws.Range("G16").Value = "'18"
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").Value = "'0.005934"
$ws.Range("E17").Value = "'0.81%"
$ws.Range("G17").Value = "'18"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "'3.384"
$ws.Range("E18").Value = "'-0.03%"
$ws.Range("G18").Value = "'18"
$ws.Range("E19").Value = "'4.12%"
$ws.Range("G19").Value = "'18"
$ws.Range("D20").Value = "'0.3328"
$ws.Range("E20").Value = "'-4.35%"
$ws.Range("G20").Value = "'18"
$ws.Range("D21").Value = "'0.1386"
$ws.Range("E21").Value = "'-2.38%"
$ws.Range("G21").Value = "'18"
$ws.Range("D22").Value = "'0.2573"
$ws.Range("E22").Value = "'2.73%"
$ws.Range("G22").Value = "'18"
$ws.Range("D23").Value = "'0.001308"
$ws.Range("E23").Value = "'4.99%"
$ws.Range("G23").Value = "'18"
$ws.Range("D24").Value = "'0.004411"
$ws.Range("E24").Value = "'0.83%"
$ws.Range("G24").Value = "'18"
$ws.Range("D25").Value = "'0.0001281"
$ws.Range("E25").Value = "'7.63%"
$ws.Range("G25").Value = "'18"
$ws.Range("D26").Value = "'0.0003741"
$ws.Range("E26").Value = "'-5.97%"
$ws.Range("G26").Value = "'18"
$ws.Range("G27").Value = "'18"
$ws.Range("G28").Value = "'18"
$ws.Range("G29").Value = "'18"
$ws.Range("G30").Value = "'18"
$ws.Range("G31").Value = "'18"
$ws.Range("G32").Value = "'18"
$ws.Range("G33").Value = "'18"
$ws.Range("G34").Value = "'18"
$ws.Range("G35").Value = "'18"
$ws.Range("G36").Value = "'18"
$ws.Range("G37").Value = "'18"
$ws.Range("D38").Value = "'0.02678"
$ws.Range("E38").Value = "'0.16%"
$ws.Range("G38").Value = "'18"
$ws.Range("D39").Value = "'0.05674"
$ws.Range("E39").Value = "'1.97%"
$ws.Range("G39").Value = "'18"
$ws.Range("D40").Value = "'0.007668"
$ws.Range("E40").Value = "'1.37%"
$ws.Range("G40").Value = "'18"
$ws.Range("D41").Value = "'0.1422"
$ws.Range("E41").Value = "'0.89%"
$ws.Range("G41").Value = "'18"
$ws.Range("D42").Value = "'0.007375"
$ws.Range("E42").Value = "'-7.13%"
$ws.Range("G42").Value = "'18"
$ws.Range("D43").Value = "'0.001986"
$ws.Range("E43").Value = "'-1.45%"
$ws.Range("G43").Value = "'18"
$ws.Range("D44").Value = "'0.008737"
$ws.Range("E44").Value = "'3.99%"
$ws.Range("G44").Value = "'18"
$ws.Range("D45").Value = "'0.00007107"
$ws.Range("E45").Value = "'-0.61%"
$ws.Range("G45").Value = "'18"
$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("E46").Value = "'0.02%"
$ws.Range("G46").Value = "'18"
$ws.Range("D47").Value = "'0.0005771"
$ws.Range("E47").Value = "'-0.69%"
$ws.Range("G47").Value = "'18"
$ws.Range("D48").Value = "'0.002521"
$ws.Range("E48").Value = "'11.05%"
$ws.Range("G48").Value = "'18"
$ws.Range("D49").Value = "'0.003416"
$ws.Range("E49").Value = "'-11.76%"
$ws.Range("G49").Value = "'18"
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("E50").Value = "'0.02%"
$ws.Range("G50").Value = "'18"
$ws.Range("D51").Value = "'0.0002001"
$ws.Range("E51").Value = "'0.02%"
$ws.Range("G51").Value = "'18"
